$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 373
$ws.Range("F3").Value = 1074
$ws.Range("F4").Value = 9573
$ws.Range("F5").Value = 205
$ws.Range("F8").Value = 6545
$ws.Range("F10").Value = 10390
$ws.Range("F11").Value = 11558
$ws.Range("F13").Value = 1194
$ws.Range("F14").Value = 5023
$ws.Range("F15").Value = 828
$ws.Range("F16").Value = 487
$ws.Range("F21").Value = 275
$ws.Range("F22").Value = 1900
$ws.Range("F23").Value = 919
$ws.Range("F24").Value = 1306
$ws.Range("F26").Value = 6
$ws.Range("F27").Value = 2070
$ws.Range("F28").Value = 443
$ws.Range("F29").Value = 656
$ws.Range("F30").Value = 2731
$ws.Range("F31").Value = 201
$ws.Range("F32").Value = 1818
$ws.Range("F34").Value = 825
$ws.Range("F35").Value = 80
$ws.Range("F36").Value = 934
$ws.Range("F37").Value = 16
$ws.Range("F38").Value = 51
$ws.Range("F39").Value = 3412
$ws.Range("F40").Value = 242
$ws.Range("F43").Value = 595
$ws.Range("F48").Value = 4235

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 30

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6068

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 373
$ws.Range("F3").Value = 1074
$ws.Range("F4").Value = 9573
$ws.Range("F9").Value = 10394
$ws.Range("F10").Value = 11558
$ws.Range("F12").Value = 1194
$ws.Range("F13").Value = 5023
$ws.Range("F14").Value = 828
$ws.Range("F15").Value = 487
$ws.Range("F18").Value = 30
$ws.Range("F21").Value = 275
$ws.Range("F22").Value = 1900
$ws.Range("F23").Value = 919
$ws.Range("F24").Value = 1306
$ws.Range("F26").Value = 2070
$ws.Range("F27").Value = 443
$ws.Range("F28").Value = 656
$ws.Range("F29").Value = 2731
$ws.Range("F30").Value = 201
$ws.Range("F31").Value = 1818
$ws.Range("F34").Value = 825
$ws.Range("F38").Value = 80
$ws.Range("F39").Value = 934
$ws.Range("F40").Value = 16
$ws.Range("F42").Value = 242
$ws.Range("F43").Value = 0
$ws.Range("F45").Value = 595
$ws.Range("F49").Value = 4235
